# Update cryptocurrency price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.968.08"
$ws.Range("E2").Value = "  -4.85%  "

$ws.Range("D3").Value = "2.207.95"
$ws.Range("E3").Value = "  -7.64%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "79.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -10.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.504"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.39%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.455"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0769"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -11.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -13.00%  "

$ws.Range("E13").Value = "  -1.86%  "

$ws.Range("D14").Value = "2.554.46"
$ws.Range("E14").Value = "  -7.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.93%  "

$ws.Range("D17").Value = "2.223.01"
$ws.Range("E17").Value = "  -7.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.708"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.34%  "

$ws.Range("D19").Value = "38.889.00"
$ws.Range("E19").Value = "  -4.85%  "

$ws.Range("D20").Value = "0.0₃0856"
$ws.Range("E20").Value = "  -6.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "223.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.19%  "

$ws.Range("E25").Value = "  +0.13%  "

$ws.Range("E26").Value = "  -10.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.99%  "

$ws.Range("E29").Value = "  -2.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "147.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.49%  "

$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.94%  "

$ws.Range("E35").Value = "  -4.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0684"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.37%  "

$ws.Range("E37").Value = "  -5.02%  "

$ws.Range("E38").Value = "  -4.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.45%  "

$ws.Range("D43").Value = "1.892.48"
$ws.Range("E43").Value = "  -4.23%  "

$ws.Range("E44").Value = "  -12.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0253"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.11%  "

$ws.Range("D49").Value = "2.420.29"
$ws.Range("E49").Value = "  -7.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.86%  "

$ws.Range("E51").Value = "  -0.21%  "
